$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "LBU_87_1",
    "LBU_87_5",
    "LBU_05_22",
    "LBU_05_08",
    "MHZ_08_02",
    "LBU_05_24",
    "MHZ_12_01",
    "LBU_05_28",
    "LBU_87_6",
    "LBU_98_4",
    "LBU_05_25",
    "LBU_96_4",
    "MHZ_12_03",
    "LBU_05_29",
    "LBU_05_11",
    "LBU_05_20",
    "MHZ_08_01",
    "LBU_05_30",
    "LBU_05_15",
    "MHZ_08_05",
    "LBU_07_01",
    "LBU_05_21",
    "LBU_05_02",
    "LBU_05_12",
    "LBU_05_27",
    "LBU_05_26",
    "MHZ_08_04",
    "LBU_96_1",
    "LBU_05_16",
    "MHZ_12_04",
    "LBU_05_18",
    "LBU_05_23",
    "LBU_05_04",
    "LBU_05_05",
    "MHZ_12_02",
    "LBU_02_3",
    "LBU_05_14",
    "LBU_01_1",
    "LBU_05_13",
    "LBU_05_03",
    "LBU_05_09",
    "LBU_05_07",
    "LBU_05_06",
    "LBU_05_10",
    "LBU_96_2",
    "LBU_87_2",
)

$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A1").Value = "hole_id"

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
